# Refresh cryptocurrency price/volume data (GitHub Actions snapshot Sat Sep 28 04:15:43 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.024.65"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "2.690.20"
$ws.Range("E3").Value = "  +1.64%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.35"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.44"
$ws.Range("E6").Value = "  +1.35%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.51%  "

$ws.Range("E9").Value = "  +5.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.03"
$ws.Range("E10").Value = "  +4.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("E11").Value = "  -1.27%  "

$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("E13").Value = "  +9.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.14"
$ws.Range("E14").Value = "  +2.38%  "

$ws.Range("D15").Value = "3.172.25"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("D16").Value = "65.879.82"
$ws.Range("E16").Value = "  +0.93%  "

$ws.Range("D17").Value = "2.687.17"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.88"
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.80"
$ws.Range("E20").Value = "  +5.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "358.72"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.33"
$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("E24").Value = "  +17.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.87"
$ws.Range("E25").Value = "  +4.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.64"
$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("E28").Value = "  +3.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.30"
$ws.Range("E29").Value = "  +0.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "538.45"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.69"
$ws.Range("E34").Value = "  +4.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  -0.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.435"
$ws.Range("E36").Value = "  +1.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.78"
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "164.71"
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("E39").Value = "  -1.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.59"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "167.98"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.17"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0634"
$ws.Range("E45").Value = "  +1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.86"
$ws.Range("E46").Value = "  +2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  +2.45%  "

$ws.Range("E48").Value = "  +1.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.657"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.85"
$ws.Range("E50").Value = "  +5.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0993"
$ws.Range("E51").Value = "  +0.85%  "
